$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.209748333333333
$ws.Range("H2").Value = 3.629245
$ws.Range("I2").Value = 0.05758082582909403
$ws.Range("J2").Value = 0.05758082582909403
$ws.Range("M2").Value = 0.2466546666666667
$ws.Range("N2").Value = 0.7399640000000001
$ws.Range("Q2").Value = 0.2983900719088889
$ws.Range("R2").Value = 2.68551064718
$ws.Range("S2").Value = 0.05758082582909403
$ws.Range("T2").Value = 0.05758082582909403

# Row 3
$ws.Range("I3").Value = 0.08883824952969446
$ws.Range("J3").Value = 0.08883824952969446
$ws.Range("M3").Value = 0.2466546666666667
$ws.Range("N3").Value = 0.7399640000000001
$ws.Range("Q3").Value = 0.4603694247822223
$ws.Range("R3").Value = 4.14332482304
$ws.Range("S3").Value = 0.08883824952969446
$ws.Range("T3").Value = 0.08883824952969446

# Row 4
$ws.Range("G4").Value = 2.562851
$ws.Range("H4").Value = 7.688553000000001
$ws.Range("I4").Value = 0.1219849393388318
$ws.Range("J4").Value = 0.1219849393388318
$ws.Range("M4").Value = 0.2466546666666667
$ws.Range("N4").Value = 0.7399640000000001
$ws.Range("Q4").Value = 0.6321391591213334
$ws.Range("R4").Value = 5.689252432092001
$ws.Range("S4").Value = 0.1219849393388318
$ws.Range("T4").Value = 0.1219849393388318

# Row 5
$ws.Range("G5").Value = 15.37051633333333
$ws.Range("H5").Value = 46.111549
$ws.Range("I5").Value = 0.7315959853023798
$ws.Range("J5").Value = 0.7315959853023797
$ws.Range("M5").Value = 0.2466546666666667
$ws.Range("N5").Value = 0.7399640000000001
$ws.Range("Q5").Value = 3.791209582692889
$ws.Range("R5").Value = 34.120886244236
$ws.Range("S5").Value = 0.7315959853023798
$ws.Range("T5").Value = 0.7315959853023797
